$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -2
$ws.Range("F15").Value = -5
$ws.Range("F22").Value = -4
$ws.Range("F24").Value = -9
$ws.Range("F39").Value = -1
$ws.Range("F42").Value = -2
$ws.Range("F48").Value = -1
$ws.Range("F52").Value = -3
$ws.Range("F54").Value = -2
$ws.Range("F55").Value = -4
$ws.Range("F58").Value = 4
$ws.Range("F60").Value = 3
$ws.Range("F62").Value = -4
$ws.Range("F68").Value = -2
$ws.Range("F70").Value = 0
$ws.Range("F71").Value = -2
$ws.Range("F76").Value = -1
$ws.Range("F79").Value = -1
